# The author's edit swaps the contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: the deck's slide master currently carries the
# "Integral" / "Red Violet" theme (stored as theme2.xml, the part the
# single slide master's relationship actually points at) and an unused
# "Office Theme" sits in theme1.xml (only referenced by the notes
# master). After the edit the active design's 12 theme colors become
# the stock "Office" palette (what used to live in theme1.xml), i.e.
# the presentation's applied look reverts from the pink/purple
# "Integral" design back to the default blue/orange "Office" design.
#
# Reproduce that with the Design/Master ColorScheme COM surface, which
# is what actually drives the live theme part (theme2.xml) backing the
# presentation's one-and-only slide master/design.

function HexToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colorScheme = $master.ColorScheme

# Target "Office Theme" color scheme, in DrawingML clrScheme slot order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
# accent6, hlink, folHlink.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $colorScheme.Colors($i + 1).RGB = HexToRgbInt $officeColors[$i]
}
